# Fixed Fatal Flaw in RGB LED Data Line, Updated BOM
# Inserts a new "Stabilizers" line item into the BOM, above the "Case" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 9 ("Case"), shifting the
# "Case" and "Total" rows down by one.
$ws.Rows("9:9").Insert() | Out-Null

# The inserted row inherits formatting from the row above (row 8);
# clear it so the new row uses the default/unstyled formatting seen
# on the other plain data rows.
$ws.Range("A9:F9").ClearFormats() | Out-Null

# Fill in the new Stabilizers line item (values entered in this order
# so new shared-string entries land in the expected order).
$ws.Range("B9").Value = "Meckeys"
$ws.Range("F9").Value = "PCB Mounted"
$ws.Range("A9").Value = "Stabilizers"
$ws.Range("E9").Value = "Stabilizers"
$ws.Range("C9").Value = 1000
$ws.Range("D9").Value = 13

# Match the saved selection state from the authored workbook.
$ws.Range("E9").Select() | Out-Null
